# Atualizacao de bases das ligas.
# For a handful of match rows, the per-match data recorded in columns
# B (id) through AC (PL_AhUnder) had been attached to the wrong fixture;
# this swaps the B:AC content between the paired rows while leaving
# column A (the sequential row index) untouched.
#
# Each column is compared before writing so that cells whose value is
# identical in both rows of a pair are left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B through AC (id .. PL_AhUnder), in sheet order.
$columns = @(
    "B","C","D","E","F","G","H","I","J","K","L","M","N","O","P",
    "Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC"
)

# Row-number pairs whose B:AC content must be exchanged.
$pairs = @(
    @(13, 14),
    @(48, 49),
    @(72, 73),
    @(74, 75),
    @(95, 96)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $columns) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")

        $val1 = $cell1.Value()
        $val2 = $cell2.Value()

        if ($val1 -ne $val2) {
            $cell1.Value = $val2
            $cell2.Value = $val1
        }
    }
}
